# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as TEXT (matches the
# original inlineStr cells), then strip the temporary text-number-format back
# to the default style so no stray "s" attribute is left on the cell.
function Set-TextValue($sheet, $ref, $val) {
    $cell = $sheet.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "65.121.45"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "3.549.30"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue $ws "D5" "598.00"
$ws.Range("E5").Value = "  -0.25%  "
Set-TextValue $ws "D6" "133.03"
$ws.Range("E6").Value = "  -5.43%  "
$ws.Range("D7").Value = "3.548.60"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("E10").Value = "  -2.58%  "
Set-TextValue $ws "D11" "7.07"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").Value = "4.153.43"
$ws.Range("E13").Value = "  -0.77%  "
Set-TextValue $ws "D14" "0.0000183"
$ws.Range("E14").Value = "  -3.31%  "
Set-TextValue $ws "D15" "26.89"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "3.552.94"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "65.238.48"
$ws.Range("E18").Value = "  -0.22%  "
Set-TextValue $ws "D19" "9.92"
$ws.Range("E19").Value = "  -4.63%  "
Set-TextValue $ws "D20" "14.36"
$ws.Range("E20").Value = "  +0.66%  "
Set-TextValue $ws "D21" "5.81"
$ws.Range("E21").Value = "  -1.16%  "
Set-TextValue $ws "D22" "390.26"
$ws.Range("E22").Value = "  -1.90%  "
Set-TextValue $ws "D23" "0.577"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("D24").Value = "3.693.78"
$ws.Range("E24").Value = "  -0.57%  "
Set-TextValue $ws "D25" "74.06"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("E26").Value = "  -0.14%  "
Set-TextValue $ws "D27" "0.0000113"
$ws.Range("E27").Value = "  -3.06%  "
Set-TextValue $ws "D28" "7.82"
$ws.Range("E28").Value = "  -0.61%  "
Set-TextValue $ws "D29" "1.56"
$ws.Range("E29").Value = "  +22.12%  "
$ws.Range("E30").Value = "  +0.07%  "
Set-TextValue $ws "D31" "8.51"
$ws.Range("E31").Value = "  +2.44%  "
Set-TextValue $ws "D32" "2.27"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "3.549.05"
$ws.Range("E33").Value = "  -1.23%  "
Set-TextValue $ws "D34" "24.00"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").Value = "  -0.03%  "
Set-TextValue $ws "D36" "0.146"
$ws.Range("E36").Value = "  -1.72%  "
Set-TextValue $ws "D37" "170.69"
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("E40").Value = "  +0.14%  "
Set-TextValue $ws "D41" "0.0812"
$ws.Range("E41").Value = "  +0.66%  "
Set-TextValue $ws "D42" "0.825"
$ws.Range("E42").Value = "  -1.11%  "
Set-TextValue $ws "D43" "26.18"
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws "D44" "1.25"
$ws.Range("E44").Value = "  +4.16%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D45" "43.03"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("D49").Value = "2.452.27"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("E51").Value = "  -0.24%  "
